$d = $word.ActiveDocument

# Locate the run of text that needs to be split into three separate runs:
#   "do Cục CSQLHC về TTXH cấp"
# becomes three runs with identical rPr (color 000000):
#   "do "  |  "Cục CSQLHC về TTXH"  |  " cấp"
$finder = $d.Content
$found = $finder.Find.Execute("do Cục CSQLHC về TTXH cấp", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Build a fresh Range object over the found span (avoid reusing a Range
    # that already has Find state attached to it).
    $target = $d.Range($finder.Start, $finder.End)

    $payload = '<?xml version="1.0" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p>' +
        '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">do </w:t></w:r>' +
        '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Cục CSQLHC về TTXH</w:t></w:r>' +
        '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> cấp</w:t></w:r>' +
        '</w:p>' +
        '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData>' +
        '</pkg:part>' +
        '</pkg:package>'

    $target.InsertXML($payload, "Replace")
}
else {
    Write-Host "Target text not found; no changes made."
}
